$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.832.75'
$ws.Range('E2').Value = '  +3.01%  '
$ws.Range('D3').Value = '3.445.02'
$ws.Range('E3').Value = '  +1.93%  '
# '1.00' would otherwise be auto-converted to a number and lose its
# trailing zero, so force the cell format to Text first.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '583.89'
$ws.Range('E5').Value = '  +2.44%  '
$ws.Range('D6').Value = '146.48'
$ws.Range('E6').Value = '  +4.37%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.476'
$ws.Range('E8').Value = '  +0.79%  '
$ws.Range('D9').Value = '7.65'
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('E11').Value = '  +2.09%  '
$ws.Range('D12').Value = '4.037.94'
$ws.Range('E12').Value = '  +2.03%  '
$ws.Range('D13').Value = '29.44'
$ws.Range('E13').Value = '  +5.80%  '
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('D15').Value = '3.449.65'
$ws.Range('E15').Value = '  +2.40%  '
$ws.Range('E16').Value = '  +2.30%  '
$ws.Range('D17').Value = '62.827.80'
$ws.Range('E17').Value = '  +2.83%  '
$ws.Range('D18').Value = '6.24'
$ws.Range('E18').Value = '  +2.61%  '
$ws.Range('E19').Value = '  +5.84%  '
$ws.Range('D20').Value = '9.31'
$ws.Range('E20').Value = '  +4.95%  '
$ws.Range('D21').Value = '394.53'
$ws.Range('E21').Value = '  +3.29%  '
$ws.Range('D22').Value = '75.37'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '0.562'
$ws.Range('E23').Value = '  +2.27%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '0.0000118'
$ws.Range('E25').Value = '  +4.06%  '
$ws.Range('D26').Value = '3.588.30'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('D28').Value = '7.73'
$ws.Range('E28').Value = '  +7.61%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = '8.19'
$ws.Range('E30').Value = '  +3.19%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '1.44'
$ws.Range('E31').Value = '  +6.62%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '2.15'
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('D34').Value = '23.79'
$ws.Range('E34').Value = '  +2.63%  '
$ws.Range('D35').Value = '5.32'
$ws.Range('E35').Value = '  +7.08%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').Value = '7.07'
$ws.Range('E36').Value = '  +2.20%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
# '1.60' would otherwise be auto-converted to a number and lose its
# trailing zero, so force the cell format to Text first.
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.60'
$ws.Range('E37').Value = '  +9.85%  '
$ws.Range('D38').Value = '168.32'
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '30.69'
$ws.Range('E39').Value = '  +18.83%  '
$ws.Range('B40').Value = 'RenzoRestakedETH'
$ws.Range('C40').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D40').Value = '3.478.01'
$ws.Range('E40').Value = '  +1.84%  '
$ws.Range('D41').Value = '0.0767'
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('D42').Value = '0.791'
$ws.Range('E42').Value = '  +1.59%  '
# '42.90' would otherwise be auto-converted to a number and lose its
# trailing zero, so force the cell format to Text first.
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.90'
$ws.Range('E43').Value = '  +1.37%  '
$ws.Range('D44').Value = '4.48'
$ws.Range('E44').Value = '  +3.20%  '
$ws.Range('E45').Value = '  +5.14%  '
$ws.Range('E46').Value = '  +7.80%  '
$ws.Range('D47').Value = '2.517.60'
$ws.Range('E47').Value = '  +3.29%  '
$ws.Range('D48').Value = '23.62'
$ws.Range('E48').Value = '  +3.72%  '
$ws.Range('D49').Value = '6.74'
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
# '1.00' would otherwise be auto-converted to a number and lose its
# trailing zero, so force the cell format to Text first.
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '0.0268'
$ws.Range('E51').Value = '  +3.00%  '
